# Updates the cryptocurrency price/volume table on Sheet1 with the latest
# scraped values (GitHub Actions refresh run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '65.215.88'
$cell.Style = "Normal"
$ws.Range("E2").Value = '  -2.08%  '

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '3.476.23'
$cell.Style = "Normal"
$ws.Range("E3").Value = '  -1.33%  '

$ws.Range("E4").Value = '  -0.06%  '

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '587.47'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  -3.28%  '

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '137.10'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  -4.62%  '

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '3.473.93'
$cell.Style = "Normal"
$ws.Range("E7").Value = '  -1.38%  '

$ws.Range("E8").Value = '  +0.08%  '

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.490'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  -3.64%  '

$ws.Range("E10").Value = '  -6.17%  '

$ws.Range("E11").Value = '  -7.14%  '

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '0.381'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  -6.09%  '

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '4.063.40'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  -1.26%  '

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '0.0000182'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  -6.87%  '

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '3.481.99'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  -1.25%  '

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '26.45'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  -7.86%  '

$ws.Range("E17").Value = '  -1.29%  '

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '65.093.61'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  -2.09%  '

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '9.71'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  -10.06%  '

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '5.77'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  -5.98%  '

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '13.88'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  -5.16%  '

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '388.56'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  -8.18%  '

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '0.553'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  -6.24%  '

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  +0.05%  '

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '72.49'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  -5.92%  '

$ws.Range("E26").Value = '  -0.20%  '

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '3.612.86'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  -1.56%  '

$ws.Range("E28").Value = '  -4.54%  '

$ws.Range("E29").Value = '  +0.12%  '

$ws.Range("E30").Value = '  -6.76%  '

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '8.19'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  -8.09%  '

$ws.Range("E32").Value = '  -10.17%  '

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '3.492.03'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  -1.06%  '

$ws.Range("E34").Value = '  -0.04%  '

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '0.144'
$cell.Style = "Normal"
$ws.Range("E35").Value = '  -6.43%  '

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '23.06'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  -4.77%  '

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '170.63'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  -1.67%  '

$ws.Range("B38").Value = 'Fetch.AI'
$ws.Range("C38").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '1.19'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  -10.55%  '

$ws.Range("B39").Value = 'Aptos'
$ws.Range("C39").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '6.83'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  -9.64%  '

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '1.46'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  -10.40%  '

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '4.75'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  -8.56%  '

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '0.0777'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  -4.10%  '

$ws.Range("E43").Value = '  -4.90%  '

$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '42.46'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  -6.70%  '

$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  -0.06%  '

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '24.81'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  +8.31%  '

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '4.36'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  -12.48%  '

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '1.62'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  -9.48%  '

$ws.Range("E49").Value = '  +2.57%  '

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '6.67'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  -5.55%  '

$ws.Range("B51").Value = 'dogwifhat'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '2.07'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  -12.72%  '

